# ADD results from server
# Updates the row-2 result values on each year sheet (2025, 2030, 2035,
# 2040, 2045, 2050) to the latest values returned by the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 3.195649908062549
$ws.Range("E2").Value = 29028.92796736782
$ws.Range("I2").Value = 26670.32387598415
$ws.Range("L2").Value = 26064.51078540457
$ws.Range("M2").Value = 11761.75092488
$ws.Range("N2").Value = 7266.369060695893
$ws.Range("O2").Value = 7167.22549913468

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 5211.62517687574
$ws.Range("E2").Value = 55931.26333373201
$ws.Range("I2").Value = 59670.73708636816
$ws.Range("L2").Value = 26064.51078540457
$ws.Range("M2").Value = 22435.25925056625
$ws.Range("N2").Value = 10913.24490122771
$ws.Range("O2").Value = 9754.483740133039

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 5669.215729593655
$ws.Range("B2").Value = 8085.2543229653
$ws.Range("E2").Value = 59953.97509265952
$ws.Range("I2").Value = 74727.52238211618
$ws.Range("L2").Value = 26064.51078540457
$ws.Range("M2").Value = 28280.46618997875
$ws.Range("N2").Value = 13303.93982533958
$ws.Range("O2").Value = 15862.05422188708

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 5669.215729593655
$ws.Range("B2").Value = 8085.2543229653
$ws.Range("E2").Value = 59953.97509265952
$ws.Range("I2").Value = 74727.52238211618
$ws.Range("L2").Value = 26064.51078540457
$ws.Range("M2").Value = 28280.46618997875
$ws.Range("N2").Value = 13303.93982533958
$ws.Range("O2").Value = 15862.05422188708

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 5669.215729593655
$ws.Range("B2").Value = 8085.2543229653
$ws.Range("E2").Value = 59953.97509265952
$ws.Range("I2").Value = 74727.52238211618
$ws.Range("L2").Value = 26064.51078540457
$ws.Range("M2").Value = 28280.46618997875
$ws.Range("N2").Value = 13303.93982533958
$ws.Range("O2").Value = 15862.05422188708

# --- Sheet "2050" ---
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 5669.215729593655
$ws.Range("B2").Value = 8085.2543229653
$ws.Range("E2").Value = 59953.97509265952
$ws.Range("I2").Value = 74727.52238211618
$ws.Range("L2").Value = 26064.51078540457
$ws.Range("M2").Value = 28280.46618997875
$ws.Range("N2").Value = 13303.93982533958
$ws.Range("O2").Value = 15862.05422188708
